$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Absent" (column H) = 1 for rows 3-11 and 13-18
foreach ($r in @(3,4,5,6,7,8,9,10,11,13,14,15,16,17,18)) {
    $ws.Cells.Item($r, 8).Value = 1
}

# Row 12 is a special case: instead of Absent, the Total Attendance Count (D)
# and Real (E) columns are marked 1
$ws.Cells.Item(12, 4).Value = 1
$ws.Cells.Item(12, 5).Value = 1
